# Apply the benchmark-results edits to the single-column results table.
# Row numbers below are 1-based table rows (each row = one cell = one value).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (Range.Text preserves the run's rPr).
$t.Cell(1, 1).Range.Text  = "0M"       # was 100
$t.Cell(2, 1).Range.Text  = "0M"       # was 0
$t.Cell(3, 1).Range.Text  = "0M"       # was 47
$t.Cell(4, 1).Range.Text  = "42"       # was 3
$t.Cell(5, 1).Range.Text  = "0.00003"  # was 0.00004
$t.Cell(6, 1).Range.Text  = "0.00009"  # was 0.00008
$t.Cell(12, 1).Range.Text = "0.00166"  # was 0.00016

# Collapse the three multi-tab detail rows down to a single summary value.
$t.Cell(44, 1).Range.Text = "100"   # was "24<tab>...<tab>100.0"
$t.Cell(45, 1).Range.Text = "0"     # was "5<tab>...<tab>100.0"
$t.Cell(46, 1).Range.Text = "47"    # was "10<tab>...<tab>100.0"
